$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" text is a plain number (e.g. "65.09") would be
# auto-parsed as a numeric value by Excel on assignment. Force those specific
# cells to Text format first so the literal text is preserved, matching the
# rest of the (always-text) Price column.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.447.07"
$ws.Range("E2").Value = "  +5.05%  "
$ws.Range("D3").Value = "2.246.45"
$ws.Range("E3").Value = "  +3.93%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "229.13"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("D7").Value = "65.09"
$ws.Range("E7").Value = "  +1.25%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +2.76%  "
$ws.Range("D10").Value = "0.0906"
$ws.Range("E10").Value = "  +5.43%  "
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "2.579.45"
$ws.Range("E12").Value = "  +3.87%  "
$ws.Range("D13").Value = "16.16"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "22.33"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").Value = "2.245.16"
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("D18").Value = "41.370.78"
$ws.Range("E18").Value = "  +5.03%  "
$ws.Range("D19").Value = "74.03"
$ws.Range("E19").Value = "  +3.10%  "
$ws.Range("D20").Value = "0.0₃0916"
$ws.Range("E20").Value = "  +7.58%  "
$ws.Range("D21").Value = "6.15"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "253.38"
$ws.Range("E22").Value = "  +9.40%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("E25").Value = "  -7.20%  "
$ws.Range("D26").Value = "9.67"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "172.82"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").Value = "0.145"
$ws.Range("E28").Value = "  +3.90%  "
$ws.Range("D29").Value = "20.49"
$ws.Range("E29").Value = "  +3.03%  "
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("D31").Value = "2.84"
$ws.Range("E31").Value = "  +6.52%  "
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("D33").Value = "4.70"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("D34").Value = "4.86"
$ws.Range("E34").Value = "  +2.06%  "
$ws.Range("D35").Value = "7.23"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("D36").Value = "0.0632"
$ws.Range("E36").Value = "  +2.22%  "
$ws.Range("E37").Value = "  +7.50%  "
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("D40").Value = "0.000234"
$ws.Range("E40").Value = "  +48.02%  "
$ws.Range("D41").Value = "4.84"
$ws.Range("E41").Value = "  +14.88%  "
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("E43").Value = "  +11.78%  "
$ws.Range("D44").Value = "17.83"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "101.91"
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("D47").Value = "1.512.30"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").Value = "0.0943"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").Value = "51.68"
$ws.Range("E51").Value = "  +11.29%  "
